$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Step 1: Insert 2 new rows before row 31. This shifts the existing
# "Status legend" rows (old 31-35) down to 33-37, and creates two new
# blank rows 31 and 32 (which will become data rows for 30-Jan / 31-Jan
# i.e. the dates 44285 / 44286).
# ------------------------------------------------------------------
$ws.Rows("31:32").Insert()

# ------------------------------------------------------------------
# Step 2: Copy row 27's current (pre-edit) formatting onto the two new
# rows 31 and 32, so they look like ordinary data rows instead of the
# placeholder row 30 formatting they inherited from the insert.
# ------------------------------------------------------------------
$ws.Range("A27:G27").Copy()
$ws.Range("A31:G31").PasteSpecial(-4122)
$ws.Range("A27:G27").Copy()
$ws.Range("A32:G32").PasteSpecial(-4122)

# ------------------------------------------------------------------
# Step 3: Row 30 was a blank placeholder row with its own special
# formatting; give it the normal data-row formatting (same as row 17).
# ------------------------------------------------------------------
$ws.Range("A17:G17").Copy()
$ws.Range("A30:G30").PasteSpecial(-4122)

# ------------------------------------------------------------------
# Step 4: Row 27's "Status" cell (F27) changes its fill style; reuse
# row 17 as a template (identical formatting elsewhere).
# ------------------------------------------------------------------
$ws.Range("A17:G17").Copy()
$ws.Range("A27:G27").PasteSpecial(-4122)

# ------------------------------------------------------------------
# Step 5: D29's border style changes; copy it from another sheet/cell
# that already uses the desired style.
# ------------------------------------------------------------------
$src = $wb.Worksheets.Item("DEC-2020")
$src.Range("D28").Copy()
$ws.Range("D29").PasteSpecial(-4122)

# ------------------------------------------------------------------
# Step 6: Fill in the actual data values.
# ------------------------------------------------------------------
$ws.Range("C27").Value = "Qmvar-2.0"
$ws.Range("D27").Value = "Dropdown menu design issues checked"
$ws.Range("E27").Value = 0.9
$ws.Range("F27").Value = "WIP"

$ws.Range("D28").Value = "Week Off"

$ws.Range("A30").Value = 29
$ws.Range("B30").Value = 44284
$ws.Range("C30").Value = "Hayaai"
$ws.Range("D30").Value = "Invoice Design created"
$ws.Range("E30").Value = 0.8
$ws.Range("F30").Value = "WIP"

$ws.Range("A31").Value = 30
$ws.Range("B31").Value = 44285

$ws.Range("A32").Value = 31
$ws.Range("B32").Value = 44286

# ------------------------------------------------------------------
# Step 7: Update the selected / active cell shown in the sheet view.
# ------------------------------------------------------------------
$ws.Range("E30").Select()

$excel.CutCopyMode = $false
